$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.140.02'
$ws.Range('E2').Value = '  +1.70%  '

$ws.Range('D3').Value = '3.270.12'
$ws.Range('E3').Value = '  +0.23%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.09'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.29%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '185.02'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.60%  '

$ws.Range('E7').Value = '  -0.06%  '

$ws.Range('E8').Value = '  -1.09%  '

$ws.Range('E9').Value = '  +3.93%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.71'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.54%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.417'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.31%  '

$ws.Range('D12').Value = '3.835.77'
$ws.Range('E12').Value = '  +0.06%  '

$ws.Range('E13').Value = '  +0.34%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.68'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.87%  '

$ws.Range('D15').Value = '68.157.49'
$ws.Range('E15').Value = '  +1.67%  '

$ws.Range('E16').Value = '  +2.81%  '

$ws.Range('D17').Value = '3.263.77'
$ws.Range('E17').Value = '  -0.07%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.87'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.10%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.63'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.19%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '383.86'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.22%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.71'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.22%  '

$ws.Range('E22').Value = '  +0.15%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.42'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.16%  '

$ws.Range('E24').Value = '  +0.34%  '

$ws.Range('E25').Value = '  +1.67%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.86'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.18%  '

$ws.Range('E27').Value = '  +2.85%  '

$ws.Range('E28').Value = '  -0.08%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.65%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.74'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.86%  '

$ws.Range('E31').Value = '  +6.48%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '22.94'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.33%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.28'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.24%  '

$ws.Range('E35').Value = '  +3.45%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '162.60'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.68%  '

$ws.Range('E37').Value = '  +0.16%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.837'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.31%  '

$ws.Range('E39').Value = '  +4.89%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.70'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.94%  '

$ws.Range('E41').Value = '  +5.31%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.62'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.61%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.52'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.40%  '

$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '348.07'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.29%  '

$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '25.48'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.45%  '

$ws.Range('D46').Value = '2.656.43'
$ws.Range('E46').Value = '  -3.56%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0688'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.55%  '

$ws.Range('E48').Value = '  +1.83%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '32.10'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.34%  '

$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.67%  '

$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.103'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.82%  '
